$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F50").Value = "Order model already contained product (MIS/CNC) and a gtt flag; we extended webhook parsing to derive product from trade_type and exposed product/gtt through the editable OrderUpdate schema and edit endpoint."
$ws.Range("G50").Value = "implemented"
$ws.Range("H50").Value = "Backoffice and APIs can now track product type and a GTT preference on orders; real GTT placement will be wired in a later sprint."
$ws.Range("I50").Value = "Implement actual Zerodha GTT order placement and tie it to the gtt flag once we integrate Kite GTT APIs."

$ws.Range("F51").Value = "Waiting Queue edit dialog now allows selecting product (MIS/CNC) and toggling a GTT preference checkbox, with changes persisted via PATCH /api/orders/{id}."
$ws.Range("G51").Value = "implemented"
$ws.Range("H51").Value = "Users can adjust trade type and mark orders as GTT-preferred before executing them from the manual queue."
$ws.Range("I51").Value = "Once real GTT execution exists, ensure the UI clearly indicates which orders will be sent as GTT vs regular/AMO."
